$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.615.43'
$ws.Range("D3").Value = '1.597.32'
$ws.Range("E3").Value = '  +0.18%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.29'
$ws.Range("E6").Value = '  +0.22%  '
$ws.Range("E8").Value = '  +0.16%  '
$ws.Range("E9").Value = '  +0.28%  '
$ws.Range("E10").Value = '  -0.82%  '
$ws.Range("E11").Value = '  +0.49%  '
$ws.Range("D12").Value = '1.821.13'
$ws.Range("E12").Value = '  +0.18%  '
$ws.Range("D13").Value = '1.582.82'
$ws.Range("E13").Value = '  -1.06%  '
$ws.Range("E14").Value = '  +0.01%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.522'
$ws.Range("E15").Value = '  -0.11%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.85'
$ws.Range("E16").Value = '  +0.14%  '
$ws.Range("D17").Value = '26.597.55'
$ws.Range("E18").Value = '  +1.19%  '
$ws.Range("E19").Value = '  +0.01%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '208.34'
$ws.Range("E20").Value = '  -0.45%  '
$ws.Range("E21").Value = '  +5.07%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.28'
$ws.Range("E22").Value = '  +0.75%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.32'
$ws.Range("E23").Value = '  +0.13%  '
$ws.Range("E24").Value = '  +0.38%  '
$ws.Range("E27").Value = '  -0.22%  '
$ws.Range("E28").Value = '  -0.36%  '
$ws.Range("E29").Value = '  -0.33%  '
$ws.Range("E30").Value = '  +1.58%  '
$ws.Range("E31").Value = '  -0.09%  '
$ws.Range("E32").Value = '  +0.57%  '
$ws.Range("E33").Value = '  +1.04%  '
$ws.Range("D34").Value = '1.283.31'
$ws.Range("E34").Value = '  -0.16%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.616'
$ws.Range("E35").Value = '  -6.86%  '
$ws.Range("E36").Value = '  +0.88%  '
$ws.Range("E37").Value = '  +0.87%  '
$ws.Range("E38").Value = '  -0.48%  '
$ws.Range("E39").Value = '  +0.67%  '
$ws.Range("E40").Value = '  +21.23%  '
$ws.Range("E41").Value = '  +2.05%  '
$ws.Range("E42").Value = '  -0.14%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '64.23'
$ws.Range("E43").Value = '  +1.20%  '
$ws.Range("E44").Value = '  -1.12%  '
$ws.Range("D45").Value = '1.733.70'
$ws.Range("E47").Value = '  -1.48%  '
$ws.Range("E48").Value = '  +4.15%  '
$ws.Range("E49").Value = '  +0.97%  '
$ws.Range("E50").Value = '  -0.09%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.42'
$ws.Range("E51").Value = '  -1.34%  '
